$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.265.47"
$ws.Range("E2").Value = "  +4.96%  "

$ws.Range("D3").Value = "2.363.11"
$ws.Range("E3").Value = "  +1.39%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.46"
$ws.Range("E5").Value = "  -0.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.22"
$ws.Range("E6").Value = "  -0.38%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  -0.54%  "

$ws.Range("E8").Value = "  -0.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.616"
$ws.Range("E9").Value = "  +0.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.26"
$ws.Range("E10").Value = "  +1.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0920"
$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.48"
$ws.Range("E12").Value = "  -0.65%  "

$ws.Range("E13").Value = "  +1.70%  "

$ws.Range("E14").Value = "  -2.37%  "

$ws.Range("D15").Value = "2.718.17"
$ws.Range("E15").Value = "  +1.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.22"
$ws.Range("E16").Value = "  -1.40%  "

$ws.Range("D17").Value = "2.359.74"
$ws.Range("E17").Value = "  +1.71%  "

$ws.Range("D18").Value = "45.165.59"
$ws.Range("E18").Value = "  +4.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.35"
$ws.Range("E19").Value = "  +8.98%  "

$ws.Range("E20").Value = "  -2.61%  "

$ws.Range("E21").Value = "  -0.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.14"
$ws.Range("E22").Value = "  -1.28%  "

$ws.Range("E23").Value = "  -0.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "259.86"
$ws.Range("E24").Value = "  -3.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.32"
$ws.Range("E25").Value = "  +2.26%  "

$ws.Range("E26").Value = "  -0.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.18"
$ws.Range("E27").Value = "  +0.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.34"
$ws.Range("E28").Value = "  -3.63%  "

$ws.Range("E29").Value = "  +2.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0964"
$ws.Range("E30").Value = "  +9.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.33"
$ws.Range("E31").Value = "  -1.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.64"
$ws.Range("E32").Value = "  -3.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "168.86"
$ws.Range("E33").Value = "  +1.03%  "

$ws.Range("E34").Value = "  +5.98%  "

$ws.Range("E35").Value = "  -0.85%  "

$ws.Range("E36").Value = "  +3.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.80"
$ws.Range("E37").Value = "  +1.76%  "

$ws.Range("E38").Value = "  +4.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.91"
$ws.Range("E39").Value = "  +4.20%  "

$ws.Range("E40").Value = "  -1.49%  "

$ws.Range("E41").Value = "  +4.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.54"
$ws.Range("E42").Value = "  -5.46%  "

$ws.Range("E43").Value = "  -1.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "69.60"
$ws.Range("E44").Value = "  -2.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.88"
$ws.Range("E45").Value = "  -3.96%  "

$ws.Range("E46").Value = "  -0.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "81.58"
$ws.Range("E47").Value = "  +5.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "112.42"
$ws.Range("E48").Value = "  -1.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.54"
$ws.Range("E49").Value = "  +4.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.25"
$ws.Range("E50").Value = "  +3.22%  "

$ws.Range("D51").Value = "1.669.44"
$ws.Range("E51").Value = "  +0.12%  "

